# Applies "hybrid bold + color" highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) inside specific bullet
# paragraphs of the resume, matching the target diff.
#
# NOTE: Write-Output/string-interpolation in this host mangles non-ASCII
# characters (bullets, plus/minus signs) purely for log display purposes;
# the underlying Range.Text / Find.Execute operations are unaffected, so
# matching logic below avoids relying on printed/interpolated copies of
# those characters and instead compares/searches the live COM strings
# directly.

$d = $word.ActiveDocument

# Word stores RGB colors as 0xBBGGRR (reverse byte order of the usual
# "RRGGBB" hex notation). Target color 2C3E50 (R=2C,G=3E,B=50) -> 0x503E2C.
$HighlightColor = 0x503E2C

# Bullet character "\u2022 " (U+2022 + space) that prefixes each bullet line.
$BulletPrefix = [string]([char]0x2022) + " "

function Strip-Bullet {
    param([string]$Text)
    $t = $Text.Trim()
    if ($t.StartsWith($BulletPrefix)) {
        $t = $t.Substring($BulletPrefix.Length)
    }
    return $t
}

# Finds the (first) paragraph whose bullet-stripped, trimmed text is
# EXACTLY $FullText, and bold+colors each substring in $Metrics in order.
function Highlight-ExactParagraph {
    param(
        [string]$FullText,
        [string[]]$Metrics
    )

    $target = $null
    foreach ($p in $d.Paragraphs) {
        $t = Strip-Bullet $p.Range.Text
        if ($t -eq $FullText) {
            $target = $p
            break
        }
    }
    if ($null -eq $target) {
        Write-Output "NOT FOUND (exact): $FullText"
        return
    }

    Apply-Metrics $target $Metrics
}

function Apply-Metrics {
    param($Paragraph, [string[]]$Metrics)

    $pEnd = $Paragraph.Range.End
    $cursor = $Paragraph.Range.Start

    foreach ($metric in $Metrics) {
        $r = $d.Range($cursor, $pEnd)
        $found = $r.Find.Execute($metric, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $r.Bold = 1
            $r.Font.Color = $HighlightColor
            $cursor = $r.End
        } else {
            Write-Output "METRIC NOT FOUND: $metric"
        }
    }
}

$PlusMinus = [string]([char]0x00B1)

# 1) "...improving demographic classification accuracy from 23% to 64%"
Highlight-ExactParagraph "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%" @("23%", "64%")

# 2) "...Achieved 87% ... 71%, reducing polling error margins from ±4.2% to ±2.1%"
Highlight-ExactParagraph ("Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from " + $PlusMinus + "4.2% to " + $PlusMinus + "2.1%") @("87%", "71%", ($PlusMinus + "4.2%"), ($PlusMinus + "2.1%"))

# 3) "Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
Highlight-ExactParagraph "Wrote RFP and analyzed bids from 1,200 vendors for research platform development" @("1,200")

# 4) "...became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+"
Highlight-ExactParagraph "Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+" @("`$400M", "`$1B")

# 5) "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
Highlight-ExactParagraph "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M" @("73.5%", "`$4.7M")

# 6) "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" (short form, no ± suffix)
Highlight-ExactParagraph "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" @("87%", "71%")
